$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (shifts old N/O/P -> O/P/Q), inheriting the
# formatting (incl. width) of the column immediately to its left, the same
# as a manual right-click "Insert" on the column header in Excel.
$mWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $mWidth

# Make "Repayment schedule" the active sheet/tab and leave the selection
# where the edit left it.
$ws.Activate() | Out-Null
$ws.Range("S5").Select() | Out-Null
